$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

function Clear-CellValue($ws, $ref) {
    $ws.Range($ref).ClearContents() | Out-Null
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
Set-CellValue $ws "H29" 2626
Set-CellValue $ws "I29" 2000
Set-CellValue $ws "J29" 4504
Set-CellValue $ws "K29" 6000
Set-CellValue $ws "L29" 13512
Set-CellValue $ws "M29" -5719
Set-CellValue $ws "N29" -14074
# Row 33
Set-CellValue $ws "H33" 367.5
Set-CellValue $ws "I33" 394.42856
Set-CellValue $ws "K33" 394.42856
Set-CellValue $ws "M33" -165.42856
# Row 38
Set-CellValue $ws "H38" 2609
Set-CellValue $ws "J38" 4999.5
Set-CellValue $ws "L38" 14998.5
Set-CellValue $ws "N38" -15742.5
# Row 43
Set-CellValue $ws "H43" 3000
Set-CellValue $ws "I43" 3000
Set-CellValue $ws "J43" 0
Set-CellValue $ws "K43" 3000
Set-CellValue $ws "L43" 0
Clear-CellValue $ws "M43"
Set-CellValue $ws "N43" -2931
# Row 53
Set-CellValue $ws "H53" 427.66666
Set-CellValue $ws "I53" 314.33334
Set-CellValue $ws "J53" 541
Set-CellValue $ws "K53" 314.33334
Set-CellValue $ws "L53" 541
Set-CellValue $ws "M53" 322.66666
Set-CellValue $ws "N53" -1815
# Row 58
Set-CellValue $ws "H58" 0
Set-CellValue $ws "J58" 0
Clear-CellValue $ws "L58"
Set-CellValue $ws "N58" 0
# Row 138
Set-CellValue $ws "H138" 4697.0586
Set-CellValue $ws "J138" 5461.5557
Set-CellValue $ws "L138" 16384.6671
Set-CellValue $ws "N138" -26664.6671
# Row 141
Set-CellValue $ws "H141" 9330
Set-CellValue $ws "I141" 7995
Set-CellValue $ws "J141" 12000
Set-CellValue $ws "K141" 23985
Set-CellValue $ws "L141" 36000
Set-CellValue $ws "M141" -18805
Set-CellValue $ws "N141" -46360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 16
Set-CellValue $ws "H16" 10000
Set-CellValue $ws "J16" 15000
Set-CellValue $ws "L16" 15000
Set-CellValue $ws "N16" -15574
# Row 36
Set-CellValue $ws "H36" 5499.5
Set-CellValue $ws "I36" 3000
Set-CellValue $ws "J36" 7999
Set-CellValue $ws "K36" 3000
Set-CellValue $ws "L36" 7999
Set-CellValue $ws "M36" -2654
Set-CellValue $ws "N36" -8691
# Row 45
Set-CellValue $ws "H45" 2319.889
Set-CellValue $ws "I45" 2148
Set-CellValue $ws "K45" 2148
Set-CellValue $ws "M45" -1771
# Row 61
Set-CellValue $ws "H61" 3636.111
Set-CellValue $ws "I61" 3590.625
Set-CellValue $ws "K61" 3590.625
Set-CellValue $ws "M61" -3378.625
# Row 63
Set-CellValue $ws "H63" 3700
Set-CellValue $ws "I63" 3700
Set-CellValue $ws "K63" 3700
Set-CellValue $ws "M63" -3014
# Row 66
Set-CellValue $ws "H66" 3700
Set-CellValue $ws "I66" 3700
Set-CellValue $ws "K66" 18500
Set-CellValue $ws "M66" -15068
# Row 74
Set-CellValue $ws "H74" 1473.4584
Set-CellValue $ws "I74" 584.05554
Set-CellValue $ws "K74" 584.05554
Set-CellValue $ws "M74" 289.94446
# Row 77
Set-CellValue $ws "H77" 1473.4584
Set-CellValue $ws "I77" 584.05554
Set-CellValue $ws "K77" 2920.2777
Set-CellValue $ws "M77" 1447.7223
# Row 102
Set-CellValue $ws "H102" 999.5
Set-CellValue $ws "J102" 999
Set-CellValue $ws "L102" 999
Set-CellValue $ws "N102" -4243
# Row 122
Set-CellValue $ws "H122" 2180.0688
Set-CellValue $ws "I122" 2193.4814
Set-CellValue $ws "J122" 1999
Set-CellValue $ws "K122" 6580.4442
Set-CellValue $ws "L122" 5997
Set-CellValue $ws "M122" -4130.4442
Set-CellValue $ws "N122" -10897
# Row 136
Set-CellValue $ws "H136" 3636.111
Set-CellValue $ws "I136" 3590.625
Set-CellValue $ws "K136" 10771.875
Set-CellValue $ws "M136" -8221.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
Set-CellValue $ws "H5" 4
Set-CellValue $ws "I5" 4
Set-CellValue $ws "K5" 4
Set-CellValue $ws "M5" 109
# Row 20
Set-CellValue $ws "H20" 2694.0908
Set-CellValue $ws "I20" 2634
Set-CellValue $ws "J20" 2799.25
Set-CellValue $ws "K20" 2634
Set-CellValue $ws "L20" 2799.25
Set-CellValue $ws "M20" -2387
Set-CellValue $ws "N20" -3293.25
# Row 82
Set-CellValue $ws "H82" 70283
Set-CellValue $ws "I82" 0
Set-CellValue $ws "K82" 0
Clear-CellValue $ws "M82"
# Row 85
Set-CellValue $ws "H85" 70283
Set-CellValue $ws "I85" 0
Set-CellValue $ws "K85" 0
Clear-CellValue $ws "M85"
# Row 131
Set-CellValue $ws "H131" 0
Set-CellValue $ws "J131" 0
Clear-CellValue $ws "L131"
Set-CellValue $ws "N131" 0
# Row 134
Set-CellValue $ws "H134" 1492.2
Set-CellValue $ws "I134" 1486.88
Set-CellValue $ws "K134" 4460.64
Set-CellValue $ws "M134" -1925.64

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
Set-CellValue $ws "H16" 6325.7144
Set-CellValue $ws "I16" 6549.8
Set-CellValue $ws "K16" 6549.8
Set-CellValue $ws "M16" -6262.8
# Row 31
Set-CellValue $ws "H31" 2734.0667
Set-CellValue $ws "I31" 2269.4614
Set-CellValue $ws "J31" 5754
Set-CellValue $ws "K31" 2269.4614
Set-CellValue $ws "L31" 5754
Set-CellValue $ws "M31" -1974.4614
Set-CellValue $ws "N31" -6344
# Row 34
Set-CellValue $ws "H34" 2734.0667
Set-CellValue $ws "I34" 2269.4614
Set-CellValue $ws "J34" 5754
Set-CellValue $ws "K34" 2269.4614
Set-CellValue $ws "L34" 5754
Set-CellValue $ws "M34" -2067.4614
Set-CellValue $ws "N34" -6158
# Row 113
Set-CellValue $ws "H113" 6325.7144
Set-CellValue $ws "I113" 6549.8
Set-CellValue $ws "K113" 6549.8
Set-CellValue $ws "M113" -4379.8
# Row 135
Set-CellValue $ws "H135" 124849.5
Set-CellValue $ws "J135" 124849.5
Set-CellValue $ws "L135" 124849.5
Set-CellValue $ws "N135" -134989.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
Set-CellValue $ws "H34" 3115.4
Set-CellValue $ws "J34" 4989
Set-CellValue $ws "L34" 14967
Set-CellValue $ws "N34" -15135
# Row 57
Set-CellValue $ws "H57" 550
Set-CellValue $ws "I57" 200
Set-CellValue $ws "J57" 900
Set-CellValue $ws "K57" 600
Set-CellValue $ws "L57" 2700
Set-CellValue $ws "M57" -41
Set-CellValue $ws "N57" -3818

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
Set-CellValue $ws "H3" 6800
Set-CellValue $ws "I3" 0
Set-CellValue $ws "J3" 6800
Set-CellValue $ws "K3" 0
Clear-CellValue $ws "L3"
Set-CellValue $ws "M3" 6800
Set-CellValue $ws "N3" -7032
# Row 97
Set-CellValue $ws "H97" 1181.2
Set-CellValue $ws "I97" 756.8889
Set-CellValue $ws "K97" 756.8889
Set-CellValue $ws "M97" -260.8889
# Row 122
Set-CellValue $ws "H122" 946.5
Set-CellValue $ws "I122" 946.5
Set-CellValue $ws "K122" 2839.5
Set-CellValue $ws "M122" -389.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 23
Set-CellValue $ws "H23" 595
Set-CellValue $ws "J23" 590
Set-CellValue $ws "L23" 590
Set-CellValue $ws "N23" -1050
# Row 43
Set-CellValue $ws "H43" 200000
Set-CellValue $ws "J43" 200000
Set-CellValue $ws "L43" 200000
Set-CellValue $ws "N43" -200386
# Row 46
Set-CellValue $ws "H46" 1999
Set-CellValue $ws "I46" 1999
Set-CellValue $ws "K46" 1999
Set-CellValue $ws "M46" -1811
# Row 132
Set-CellValue $ws "I132" 4889
Set-CellValue $ws "J132" 5998.6665
Set-CellValue $ws "K132" 14667
Set-CellValue $ws "L132" 17995.9995
Set-CellValue $ws "M132" -12137
Set-CellValue $ws "N132" -23055.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
Set-CellValue $ws "H81" 9224.637000000001
Set-CellValue $ws "I81" 6558.375
Set-CellValue $ws "K81" 13116.75
Set-CellValue $ws "M81" -12055.75
# Row 84
Set-CellValue $ws "H84" 9224.637000000001
Set-CellValue $ws "I84" 6558.375
Set-CellValue $ws "K84" 65583.75
Set-CellValue $ws "M84" -60279.75
# Row 113
Set-CellValue $ws "H113" 902
Set-CellValue $ws "I113" 824.8333
Set-CellValue $ws "J113" 1056.3334
Set-CellValue $ws "K113" 2474.4999
Set-CellValue $ws "L113" 3169.0002
Set-CellValue $ws "M113" -304.4998999999998
Set-CellValue $ws "N113" -7509.0002
# Row 122
Set-CellValue $ws "H122" 2636.6365
Set-CellValue $ws "I122" 2636.6365
Set-CellValue $ws "K122" 7909.9095
Set-CellValue $ws "M122" -5459.9095
# Row 126
Set-CellValue $ws "H126" 2050.5715
Set-CellValue $ws "I126" 2050.5715
Set-CellValue $ws "K126" 6151.7145
Set-CellValue $ws "M126" -3681.7145
# Row 136
Set-CellValue $ws "H136" 10118.941
Set-CellValue $ws "I136" 10546.357
Set-CellValue $ws "K136" 31639.071
Set-CellValue $ws "M136" -29089.071
